$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Задание 1")
$ws2 = $wb.Worksheets.Item("Задание 2")

# ---------------------------------------------------------------------
# Sheet "Задание 2": insert new column K "Ki*Vi" = Ci*Bi for both tables
# ---------------------------------------------------------------------

# Table 1 header (row 1) + data (rows 2:6)
$ws2.Range("K1").VerticalAlignment   = -4108   # xlCenter
$ws2.Range("K1").HorizontalAlignment = -4108   # xlCenter
$ws2.Range("K1").Value = "Ki*Vi"

$ws2.Range("K2").Formula  = '=C2*B2'
$ws2.Range("K3").Formula  = '=C3*B3'
$ws2.Range("K4").Formula  = '=C4*B4'
$ws2.Range("K5").Formula  = '=C5*B5'
$ws2.Range("K6").Formula  = '=C6*B6'

# Table 2 header (row 8) + data (rows 9:13)
$ws2.Range("K8").VerticalAlignment   = -4108   # xlCenter
$ws2.Range("K8").HorizontalAlignment = -4108   # xlCenter
$ws2.Range("K8").Value = "Ki*Vi"

$ws2.Range("K9").Formula  = '=C9*B9'
$ws2.Range("K10").Formula = '=C10*B10'
$ws2.Range("K11").Formula = '=C11*B11'
$ws2.Range("K12").Formula = '=C12*B12'
$ws2.Range("K13").Formula = '=C13*B13'

# ---------------------------------------------------------------------
# Sheet "Задание 1": L21 literal 0.2 -> formula referencing K21,J21,R39
# ---------------------------------------------------------------------
$ws1.Range("L21").Formula = '=($K$21-$J$21)/$R$39'

# ---------------------------------------------------------------------
# View state: active sheet becomes "Задание 1", selections updated
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("C26").Select()

$ws1.Activate()
$ws1.Range("P14").Select()
